$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290, shifting existing rows 290.. down by one.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new record's data.
$ws.Cells.Item(290, 1).Value = 6
$ws.Cells.Item(290, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(290, 3).Value = "Metropolitana"
$ws.Cells.Item(290, 4).Value = 44826
$ws.Cells.Item(290, 5).Value = 13
$ws.Cells.Item(290, 6).Value = 100112026
$ws.Cells.Item(290, 7).Value = "Haba"
$ws.Cells.Item(290, 8).Value = "Sin especificar"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 1150
$ws.Cells.Item(290, 11).Value = 8000
$ws.Cells.Item(290, 12).Value = 9000
$ws.Cells.Item(290, 13).Value = 8417
$ws.Cells.Item(290, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(290, 15).Value = "Región Metropolitana"
$ws.Cells.Item(290, 16).Value = 337
$ws.Cells.Item(290, 17).Value = 25
$ws.Cells.Item(290, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date-time number format as the rest of column D.
$ws.Cells.Item(290, 4).NumberFormat = $ws.Cells.Item(291, 4).NumberFormat()
